$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume snapshot refresh (GitHub Actions bot).
# Column D ("Price") holds numeric-looking text (e.g. thousand-dot formatted
# prices); re-applying "@" (Text) format right before the write keeps Excel
# from auto-converting it to a real number, then ClearFormats drops the extra
# style index again so the cell matches its original (unstyled) appearance.

$ws.Range("D2").Value = "30.631.51"
$ws.Range("E2").Value = "  +2.25%  "
$ws.Range("D3").Value = "1.893.71"
$ws.Range("E3").Value = "  +1.06%  "
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "244.23"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.22%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("E7").Value = "  +0.34%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2960"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.92%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06825"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +3.90%  "
$ws.Range("D10").Value = "1.892.11"
$ws.Range("E10").Value = "  +1.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "17.06"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.33%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07335"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +2.37%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "91.30"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +6.29%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.102"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.32%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6744"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +2.55%  "
$ws.Range("D16").Value = "30.632.42"
$ws.Range("E16").Value = "  +2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.000007908"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.97%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E19").Value = "  +4.92%  "
$ws.Range("D20").Value = "2.135.11"
$ws.Range("E20").Value = "  +0.95%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.872"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +3.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "178.19"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +32.41%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.072"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +9.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.295"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +2.98%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "154.22"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +2.98%  "
$ws.Range("E27").Value = "  +12.81%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.927"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.70%  "
$ws.Range("E29").Value = "  +1.25%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.340"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +4.77%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08947"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.035"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +2.83%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05206"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.7400"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +6.09%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.135"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +4.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.673"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.69%  "
$ws.Range("E37").Value = "  +10.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.703"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.63%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.169"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.39%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9335"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.86%  "
$ws.Range("E41").Value = "  +4.66%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "106.04"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +4.56%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.804"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -2.23%  "
$ws.Range("E44").Value = "  +0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.656"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.82%  "
$ws.Range("E46").Value = "  +8.15%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05845"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.60%  "
$ws.Range("E48").Value = "  +3.44%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.3891"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.485"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +5.20%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.381"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +3.92%  "
